$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Delete the final paragraph that held the second "Video" link
#    (hyperlink to https://youtu.be/06tQzrrvH0k, the " " run and "Video 2").
#    That paragraph is the 7th paragraph in the document before any other
#    edits, so remove it first while the indices are still predictable.
# ---------------------------------------------------------------------------
$d.Paragraphs(7).Range.Delete()

# ---------------------------------------------------------------------------
# 2. Trim the trailing "1" off of "Video 1" so it reads "Video " (the space
#    is preserved, matching the target xml:space="preserve" run).
# ---------------------------------------------------------------------------
$videoPara = $d.Paragraphs(6).Range
$videoFind = $videoPara.Duplicate
$videoFind.Find.Execute("Video 1") | Out-Null
$trailingDigit = $d.Range($videoFind.End - 1, $videoFind.End)
$trailingDigit.Delete()

# ---------------------------------------------------------------------------
# 3. Re-point the third hyperlink (https://www.youtube.com/watch?v=x9dE20pXdo4)
#    at the new short link, updating both the stored target and the visible
#    display text.
# ---------------------------------------------------------------------------
$videoLink = $d.Hyperlinks(3)
$videoLink.Address = "https://youtu.be/0TPt67EQbhg"

$urlFind = $d.Paragraphs(6).Range.Duplicate
$urlFind.Find.Execute(
    "https://www.youtube.com/watch?v=x9dE20pXdo4",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "https://youtu.be/0TPt67EQbhg", 2) | Out-Null
